$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 5.418
$ws.Range("D3").Value = -7.255999999999998
$ws.Range("B21").Value = 9.458
$ws.Range("B23").Value = 7.891999999999999
$ws.Range("D24").Value = -7.571
$ws.Range("B25").Value = 6.396
$ws.Range("C27").Value = -13.055
$ws.Range("C31").Value = -13.379
$ws.Range("C39").Value = -12.847
$ws.Range("C48").Value = -11.216
$ws.Range("C51").Value = -11.329
$ws.Range("C52").Value = -11.389
$ws.Range("B53").Value = 6.139999999999999
$ws.Range("C55").Value = -13.368
$ws.Range("C56").Value = -12.996
$ws.Range("B57").Value = 5.178
$ws.Range("C57").Value = -13.716
$ws.Range("D57").Value = -8.323
$ws.Range("B59").Value = 4.714
$ws.Range("D61").Value = -7.739
$ws.Range("B69").Value = 5.667000000000001
$ws.Range("D70").Value = -7.188000000000001
$ws.Range("C73").Value = -12.818
$ws.Range("B79").Value = 5.760000000000001
$ws.Range("B83").Value = 5.702
$ws.Range("D86").Value = -8.241
$ws.Range("C89").Value = -10.953
$ws.Range("C90").Value = -12.91
$ws.Range("B93").Value = 5.659000000000001
$ws.Range("D98").Value = -8.397
$ws.Range("D100").Value = -8.361999999999998
$ws.Range("D102").Value = -7.805000000000001
